# Update Leve profit-tracking figures across several sheets in the
# "Marilith_Profits" workbook (scheduled runner refresh of market-board
# averages and computed Leve prices/profits).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1347.6666
$ws.Range("I15").Value = 1347.6666
$ws.Range("K15").Value = 4042.9998
$ws.Range("M15").Value = -3873.9998

# Row 51
$ws.Range("H51").Value = 333
$ws.Range("I51").Value = 333
$ws.Range("K51").Value = 333
$ws.Range("M51").Value = 151

# Row 62
$ws.Range("H62").Value = 3527.7144
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("N62").Value = -5248

# Row 65
$ws.Range("H65").Value = 3527.7144
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 20000
$ws.Range("N65").Value = -26240

# Row 70
$ws.Range("H70").Value = 3299.6667
$ws.Range("I70").Value = 2270.8572
$ws.Range("K70").Value = 6812.571599999999
$ws.Range("M70").Value = -6542.571599999999

# Row 73
$ws.Range("H73").Value = 3299.6667
$ws.Range("I73").Value = 2270.8572
$ws.Range("K73").Value = 6812.571599999999
$ws.Range("M73").Value = -5876.571599999999

# Row 98
$ws.Range("H98").Value = 1314.15
$ws.Range("I98").Value = 1114.8948
$ws.Range("K98").Value = 1114.8948
$ws.Range("M98").Value = 383.1052

# Row 107
$ws.Range("H107").Value = 2005.6666
$ws.Range("I107").Value = 823.625
$ws.Range("J107").Value = 3356.5715
$ws.Range("K107").Value = 823.625
$ws.Range("L107").Value = 3356.5715
$ws.Range("M107").Value = 1096.375
$ws.Range("N107").Value = -7196.5715

# Row 122
$ws.Range("H122").Value = 1314.15
$ws.Range("I122").Value = 1114.8948
$ws.Range("K122").Value = 3344.6844
$ws.Range("M122").Value = -894.6844000000001

# Row 132
$ws.Range("H132").Value = 2250.6667
$ws.Range("I132").Value = 2360.8
$ws.Range("K132").Value = 7082.400000000001
$ws.Range("M132").Value = -4552.400000000001

# Row 135
$ws.Range("H135").Value = 696.4
$ws.Range("I135").Value = 601
$ws.Range("J135").Value = 839.5
$ws.Range("K135").Value = 5409
$ws.Range("L135").Value = 7555.5
$ws.Range("M135").Value = -2874
$ws.Range("N135").Value = -12625.5

# Row 141
$ws.Range("H141").Value = 2492.24
$ws.Range("I141").Value = 2179.4167
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 6538.250100000001
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -1358.250100000001
$ws.Range("N141").Value = -40360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 1624.75
$ws.Range("I63").Value = 1235.5555
$ws.Range("J63").Value = 2792.3333
$ws.Range("K63").Value = 1235.5555
$ws.Range("L63").Value = 2792.3333
$ws.Range("M63").Value = -549.5554999999999
$ws.Range("N63").Value = -4164.3333

# Row 66
$ws.Range("H66").Value = 1624.75
$ws.Range("I66").Value = 1235.5555
$ws.Range("J66").Value = 2792.3333
$ws.Range("K66").Value = 6177.7775
$ws.Range("L66").Value = 13961.6665
$ws.Range("M66").Value = -2745.7775
$ws.Range("N66").Value = -20825.6665

# Row 88
$ws.Range("H88").Value = 2371.2354
$ws.Range("I88").Value = 900.625
$ws.Range("K88").Value = 900.625
$ws.Range("M88").Value = -494.625

# Row 91
$ws.Range("H91").Value = 2371.2354
$ws.Range("I91").Value = 900.625
$ws.Range("K91").Value = 900.625
$ws.Range("M91").Value = 503.375

# Row 117
$ws.Range("H117").Value = 21165.334
$ws.Range("J117").Value = 21165.334
$ws.Range("L117").Value = 21165.334
$ws.Range("N117").Value = -30343.334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 18864
$ws.Range("J28").Value = 18864
$ws.Range("L28").Value = 18864
$ws.Range("N28").Value = -19354

# Row 31
$ws.Range("H31").Value = 2957.1614
$ws.Range("I31").Value = 2158.2222
$ws.Range("J31").Value = 8350
$ws.Range("K31").Value = 2158.2222
$ws.Range("L31").Value = 8350
$ws.Range("M31").Value = -1863.2222
$ws.Range("N31").Value = -8940

# Row 34
$ws.Range("H34").Value = 2957.1614
$ws.Range("I34").Value = 2158.2222
$ws.Range("J34").Value = 8350
$ws.Range("K34").Value = 2158.2222
$ws.Range("L34").Value = 8350
$ws.Range("M34").Value = -1956.2222
$ws.Range("N34").Value = -8754

# Row 87
$ws.Range("H87").Value = 47500
$ws.Range("J87").Value = 47500
$ws.Range("L87").Value = 47500
$ws.Range("N87").Value = -49872

# Row 90
$ws.Range("H90").Value = 47500
$ws.Range("J90").Value = 47500
$ws.Range("L90").Value = 142500
$ws.Range("N90").Value = -154356

# Row 133
$ws.Range("H133").Value = 47888
$ws.Range("J133").Value = 47888
$ws.Range("L133").Value = 47888
$ws.Range("N133").Value = -52948

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 5118.1665
$ws.Range("I75").Value = 4799.5
$ws.Range("J75").Value = 5277.5
$ws.Range("K75").Value = 14398.5
$ws.Range("L75").Value = 15832.5
$ws.Range("M75").Value = -13400.5
$ws.Range("N75").Value = -17828.5

# Row 78
$ws.Range("H78").Value = 5118.1665
$ws.Range("I78").Value = 4799.5
$ws.Range("J78").Value = 5277.5
$ws.Range("K78").Value = 43195.5
$ws.Range("L78").Value = 47497.5
$ws.Range("M78").Value = -38203.5
$ws.Range("N78").Value = -57481.5

# Row 132
$ws.Range("H132").Value = 1042.7778
$ws.Range("I132").Value = 494.5
$ws.Range("J132").Value = 1199.4286
$ws.Range("K132").Value = 4450.5
$ws.Range("L132").Value = 10794.8574
$ws.Range("M132").Value = -1920.5
$ws.Range("N132").Value = -15854.8574

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 3665.3333
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 3665.3333
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3665.3333
$ws.Range("N12").Value = -4005.3333
$ws.Range("M12").ClearContents()

# Row 22
$ws.Range("H22").Value = 2399.0588
$ws.Range("I22").Value = 2541.5715
$ws.Range("J22").Value = 1734
$ws.Range("K22").Value = 2541.5715
$ws.Range("L22").Value = 1734
$ws.Range("M22").Value = -2246.5715
$ws.Range("N22").Value = -2324

# Row 27
$ws.Range("H27").Value = 2399.0588
$ws.Range("I27").Value = 2541.5715
$ws.Range("J27").Value = 1734
$ws.Range("K27").Value = 2541.5715
$ws.Range("L27").Value = 1734
$ws.Range("M27").Value = -2434.5715
$ws.Range("N27").Value = -1948

# Row 40
$ws.Range("H40").Value = 3125
$ws.Range("I40").Value = 2665.3333
$ws.Range("J40").Value = 4504
$ws.Range("K40").Value = 2665.3333
$ws.Range("L40").Value = 4504
$ws.Range("M40").Value = -2529.3333
$ws.Range("N40").Value = -4776

# Row 46
$ws.Range("H46").Value = 3517.5
$ws.Range("I46").Value = 2646.6667
$ws.Range("J46").Value = 4205
$ws.Range("K46").Value = 2646.6667
$ws.Range("L46").Value = 4205
$ws.Range("M46").Value = -2458.6667
$ws.Range("N46").Value = -4581

# Row 93
$ws.Range("H93").Value = 1500
$ws.Range("J93").Value = 1500
$ws.Range("L93").Value = 1500
$ws.Range("N93").Value = -3996

# Row 115
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -42350

# Row 122
$ws.Range("H122").Value = 3645.8572
$ws.Range("I122").Value = 3302.4
$ws.Range("J122").Value = 4504.5
$ws.Range("K122").Value = 9907.2
$ws.Range("L122").Value = 13513.5
$ws.Range("M122").Value = -7457.200000000001
$ws.Range("N122").Value = -18413.5

# Row 132
$ws.Range("H132").Value = 8166.6665
$ws.Range("I132").Value = 7875
$ws.Range("J132").Value = 8400
$ws.Range("K132").Value = 23625
$ws.Range("L132").Value = 25200
$ws.Range("M132").Value = -21095
$ws.Range("N132").Value = -30260
